$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Case_2_10 (380 kV) vm_pu results: columns B-F and I-N, rows 2-25.
# Column A (bus index) and column G stay unchanged per the commit diff.
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035940116698993
$ws.Range("D2").Value = 1.042567694981916
$ws.Range("E2").Value = 1.051155885153868
$ws.Range("F2").Value = 1.05621891788469
$ws.Range("I2").Value = 1.034260130907798
$ws.Range("J2").Value = 1.041051150281763
$ws.Range("K2").Value = 1.045344016534309
$ws.Range("L2").Value = 1.053908169894613
$ws.Range("M2").Value = 1.058957230546133
$ws.Range("N2").Value = 1.017568119792965

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037076386633553
$ws.Range("D3").Value = 1.043421473233834
$ws.Range("E3").Value = 1.0521981059167
$ws.Range("F3").Value = 1.057275462344669
$ws.Range("I3").Value = 1.034446253596193
$ws.Range("J3").Value = 1.041830235974473
$ws.Range("K3").Value = 1.046008540291546
$ws.Range("L3").Value = 1.054762389668289
$ws.Range("M3").Value = 1.059826752007912
$ws.Range("N3").Value = 1.01783181130236

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.037811621130496
$ws.Range("D4").Value = 1.043973777974765
$ws.Range("E4").Value = 1.052873305296313
$ws.Range("F4").Value = 1.057959624398905
$ws.Range("I4").Value = 1.034565244654871
$ws.Range("J4").Value = 1.042333817253804
$ws.Range("K4").Value = 1.046437743653282
$ws.Range("L4").Value = 1.055315334915582
$ws.Range("M4").Value = 1.060389295138784
$ws.Range("N4").Value = 1.018002118579506

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038120712579849
$ws.Range("D5").Value = 1.044205931450816
$ws.Range("E5").Value = 1.053157353511057
$ws.Range("F5").Value = 1.058247367422428
$ws.Range("I5").Value = 1.034614922938505
$ws.Range("J5").Value = 1.042545393874623
$ws.Range("K5").Value = 1.046617992014419
$ws.Range("L5").Value = 1.055547842732119
$ws.Range("M5").Value = 1.060625765339699
$ws.Range("N5").Value = 1.018073639403056

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038172610351532
$ws.Range("D6").Value = 1.044244908954938
$ws.Range("E6").Value = 1.053205057849752
$ws.Range("F6").Value = 1.058295687851763
$ws.Range("I6").Value = 1.034623243872729
$ws.Range("J6").Value = 1.042580910963008
$ws.Range("K6").Value = 1.046648245437227
$ws.Range("L6").Value = 1.055586884723998
$ws.Range("M6").Value = 1.060665468379568
$ws.Range("N6").Value = 1.018085643585225

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037815751228998
$ws.Range("D7").Value = 1.043976880160468
$ws.Range("E7").Value = 1.052877099998872
$ws.Range("F7").Value = 1.057963468758802
$ws.Range("I7").Value = 1.034565909816529
$ws.Range("J7").Value = 1.042336644857462
$ws.Range("K7").Value = 1.046440152882421
$ws.Range("L7").Value = 1.055318441501559
$ws.Range("M7").Value = 1.060392454955406
$ws.Range("N7").Value = 1.018003074544332

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03632412661599
$ws.Range("D8").Value = 1.042856263587833
$ws.Range("E8").Value = 1.051507940111149
$ws.Range("F8").Value = 1.056575876956641
$ws.Range("I8").Value = 1.034323330582313
$ws.Range("J8").Value = 1.041314557687761
$ws.Range("K8").Value = 1.04556875854101
$ws.Range("L8").Value = 1.054196814325445
$ws.Range("M8").Value = 1.05925110909259
$ws.Range("N8").Value = 1.017657301549137

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.033695594463106
$ws.Range("D9").Value = 1.040880480652246
$ws.Range("E9").Value = 1.049101539480903
$ws.Range("F9").Value = 1.05413465150334
$ws.Range("I9").Value = 1.033884831812935
$ws.Range("J9").Value = 1.039509372186849
$ws.Range("K9").Value = 1.04402721943913
$ws.Range("L9").Value = 1.052221962436767
$ws.Range("M9").Value = 1.057239184816775
$ws.Range("N9").Value = 1.017045565724795

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.031943119081516
$ws.Range("D10").Value = 1.03956255641367
$ws.Range("E10").Value = 1.047501472567248
$ws.Range("F10").Value = 1.052509785366022
$ws.Range("I10").Value = 1.033585082508328
$ws.Range("J10").Value = 1.038303120781369
$ws.Range("K10").Value = 1.042995474497385
$ws.Range("L10").Value = 1.050906474842514
$ws.Range("M10").Value = 1.055897415253202
$ws.Range("N10").Value = 1.016636102880346

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.031184236407044
$ws.Range("D11").Value = 1.038991708062522
$ws.Range("E11").Value = 1.046809623118479
$ws.Range("F11").Value = 1.051806820875276
$ws.Range("I11").Value = 1.033453529924203
$ws.Range("J11").Value = 1.037780134162596
$ws.Range("K11").Value = 1.042547756175901
$ws.Range("L11").Value = 1.050337110805088
$ws.Range("M11").Value = 1.055316298326325
$ws.Range("N11").Value = 1.01645841215788

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030902344974453
$ws.Range("D12").Value = 1.038779642804382
$ws.Range("E12").Value = 1.046552788156532
$ws.Range("F12").Value = 1.051545800558583
$ws.Range("I12").Value = 1.0334044011938
$ws.Range("J12").Value = 1.037585772300899
$ws.Range("K12").Value = 1.042381308409977
$ws.Range("L12").Value = 1.050125661261288
$ws.Range("M12").Value = 1.055100427173698
$ws.Range("N12").Value = 1.01639235118153

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.030962812053147
$ws.Range("D13").Value = 1.038825132732489
$ws.Range("E13").Value = 1.046607873397806
$ws.Range("F13").Value = 1.051601786141906
$ws.Range("I13").Value = 1.033414951430003
$ws.Range("J13").Value = 1.037627468172783
$ws.Range("K13").Value = 1.042417018595497
$ws.Range("L13").Value = 1.050171016205304
$ws.Range("M13").Value = 1.055146733106045
$ws.Range("N13").Value = 1.016406524143624

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.031160935356123
$ws.Range("D14").Value = 1.03897417922783
$ws.Range("E14").Value = 1.046788390026522
$ws.Range("F14").Value = 1.051785242978633
$ws.Range("I14").Value = 1.033449474318214
$ws.Range("J14").Value = 1.037764070220298
$ws.Range("K14").Value = 1.042534000515352
$ws.Range("L14").Value = 1.050319631552967
$ws.Range("M14").Value = 1.055298454726759
$ws.Range("N14").Value = 1.016452952735316

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.031283004483659
$ws.Range("D15").Value = 1.039066008147732
$ws.Range("E15").Value = 1.04689963198654
$ws.Range("F15").Value = 1.051898288976508
$ws.Range("I15").Value = 1.033470709995745
$ws.Range("J15").Value = 1.037848221808518
$ws.Range("K15").Value = 1.042606057686768
$ws.Range("L15").Value = 1.050411203355416
$ws.Range("M15").Value = 1.055391932975346
$ws.Range("N15").Value = 1.016481551139269

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.031993482483443
$ws.Range("D16").Value = 1.039600437982917
$ws.Range("E16").Value = 1.047547409194753
$ws.Range("F16").Value = 1.052556451712877
$ws.Range("I16").Value = 1.033593776160942
$ws.Range("J16").Value = 1.038337815491367
$ws.Range("K16").Value = 1.043025167725853
$ws.Range("L16").Value = 1.050944266936453
$ws.Range("M16").Value = 1.055935979490147
$ws.Range("N16").Value = 1.016647887394821

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032439132513083
$ws.Range("D17").Value = 1.039935623955442
$ws.Range("E17").Value = 1.04795400778245
$ws.Range("F17").Value = 1.052969463898708
$ws.Range("I17").Value = 1.033670501374206
$ws.Range("J17").Value = 1.038644744521211
$ws.Range("K17").Value = 1.04328780568125
$ws.Range("L17").Value = 1.051278710730634
$ws.Range("M17").Value = 1.056277212770374
$ws.Range("N17").Value = 1.01675212112438

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.032699067963281
$ws.Range("D18").Value = 1.040131115040565
$ws.Range("E18").Value = 1.048191265389713
$ws.Range("F18").Value = 1.053210426100261
$ws.Range("I18").Value = 1.03371508415649
$ws.Range("J18").Value = 1.038823706240411
$ws.Range("K18").Value = 1.043440904765002
$ws.Range("L18").Value = 1.051473810337668
$ws.Range("M18").Value = 1.056476236622382
$ws.Range("N18").Value = 1.016812881210386

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.032787698478958
$ws.Range("D19").Value = 1.040197769526042
$ws.Range("E19").Value = 1.048272180306703
$ws.Range("F19").Value = 1.05329259802643
$ws.Range("I19").Value = 1.03373025694134
$ws.Range("J19").Value = 1.038884716584797
$ws.Range("K19").Value = 1.043493091792848
$ws.Range("L19").Value = 1.051540338368589
$ws.Range("M19").Value = 1.056544096624558
$ws.Range("N19").Value = 1.016833592438888

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.032391318974345
$ws.Range("D20").Value = 1.039899663453619
$ws.Range("E20").Value = 1.047910373720901
$ws.Range("F20").Value = 1.052925145479381
$ws.Range("I20").Value = 1.03366228704118
$ws.Range("J20").Value = 1.038611820645563
$ws.Range("K20").Value = 1.043259636754294
$ws.Range("L20").Value = 1.05124282557388
$ws.Range("M20").Value = 1.05624060288479
$ws.Range("N20").Value = 1.016740941729874

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03110259322915
$ws.Range("D21").Value = 1.038930289488269
$ws.Range("E21").Value = 1.046735228273175
$ws.Range("F21").Value = 1.0517312169744
$ws.Range("I21").Value = 1.033439315482101
$ws.Range("J21").Value = 1.037723847122991
$ws.Range("K21").Value = 1.042499556266167
$ws.Range("L21").Value = 1.050275866999505
$ws.Range("M21").Value = 1.055253777001546
$ws.Range("N21").Value = 1.016439282293883

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.030292269741264
$ws.Range("D22").Value = 1.038320650778464
$ws.Range("E22").Value = 1.045997228172174
$ws.Range("F22").Value = 1.05098107891842
$ws.Range("I22").Value = 1.033297595607934
$ws.Range("J22").Value = 1.037164956365922
$ws.Range("K22").Value = 1.042020823081491
$ws.Range("L22").Value = 1.049668119336526
$ws.Range("M22").Value = 1.05463321357343
$ws.Range("N22").Value = 1.016249277214634

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.030721842743163
$ws.Range("D23").Value = 1.038643846453485
$ws.Range("E23").Value = 1.046388374432385
$ws.Range("F23").Value = 1.051378690897968
$ws.Range("I23").Value = 1.033372868890998
$ws.Range("J23").Value = 1.037461290668807
$ws.Range("K23").Value = 1.042274688303705
$ws.Range("L23").Value = 1.049990277270234
$ws.Range("M23").Value = 1.054962196162845
$ws.Range("N23").Value = 1.016350034728639

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.032412923868442
$ws.Range("D24").Value = 1.039915912509573
$ws.Range("E24").Value = 1.047930089779799
$ws.Range("F24").Value = 1.052945170882379
$ws.Range("I24").Value = 1.033665999269097
$ws.Range("J24").Value = 1.038626697727761
$ws.Range("K24").Value = 1.043272365367667
$ws.Range("L24").Value = 1.051259040456649
$ws.Range("M24").Value = 1.056257145352365
$ws.Range("N24").Value = 1.016745993333424

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.034375150331085
$ws.Range("D25").Value = 1.041391398493597
$ws.Range("E25").Value = 1.049722912054843
$ws.Range("F25").Value = 1.054765305664125
$ws.Range("I25").Value = 1.033999502100876
$ws.Range("J25").Value = 1.039976547505839
$ws.Range("K25").Value = 1.044426458663698
$ws.Range("L25").Value = 1.052732318977308
$ws.Range("M25").Value = 1.057759401403743
$ws.Range("N25").Value = 1.017204003068015
